$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.510.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.193.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.61%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.53'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.08%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.562'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.190.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.53%  '

$ws.Range("E10").Value = '  +2.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.92'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.98%  '

$ws.Range("E12").Value = '  +4.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000268'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.714.33'
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.507.73'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.82%  '

$ws.Range("E17").Value = '  +5.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.194.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.69%  '

$ws.Range("E19").Value = '  +0.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '517.21'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.97%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.739'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.93%  '

$ws.Range("E28").Value = '  +4.02%  '

$ws.Range("E29").Value = '  +9.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +15.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.48%  '

$ws.Range("E33").Value = '  +2.65%  '

$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '511.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.25%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '55.02'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0906'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.16%  '

$ws.Range("E39").Value = '  +2.94%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.128'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.88%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.90'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.302'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.86%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0670'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +16.00%  '

$ws.Range("E45").Value = '  +2.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.909.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.97%  '

$ws.Range("E47").Value = '  +3.32%  '

$ws.Range("E48").Value = '  +3.94%  '

$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.53%  '

$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.67'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.27%  '
